# Updates cryptos list values (Price column D, Volume(1h) column E)
# per commit "Updated cryptos list on Sun Aug 25 05:10:40 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.087.26"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.759.01"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -13.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("D13").Value = "3.248.25"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "63.761.76"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "2.762.99"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "361.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.550"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.99%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.0₃0931"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.993"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "331.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0258"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +0.67%  "
